$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.981.05"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.647.08"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.881.39"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "1.660.22"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.58"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "27.987.62"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.99"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").Value = "1.444.92"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.889"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.560"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.919"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.45"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("D48").Value = "1.789.58"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("E51").Value = "  +0.15%  "
